# Add header row (Title / Author / Price) above the existing book data,
# and make the new headers bold, 16pt Calibri.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Title"
$ws.Range("C1").Value = "Author"
$ws.Range("D1").Value = "Price"

$headerRange = $ws.Range("B1:D1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 16
